# Generate Report for Handback
# Updates the localization-status workbook after a failed handback transform:
#  - Status for the e50d8f56... row flips from "Ready for handoff" to
#    "Handback transform failed" on both the zh-cn and de-de sheets
#    (and, as a consequence of shared-string reuse, on the Overview sheet too).
#  - The Error Detail column (P) for that same row gets a diagnostic message,
#    one per locale.
#  - Error Detail column width is widened to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status column (C) for the e50d8f56 row (row 3) on both locale sheets.
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# The Overview sheet's zh-cn/de-de status columns (E/F) for the same row
# shared this exact string, so they pick up the same new status text.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# Error Detail column (P) for the e50d8f56 row (row 3) on both locale sheets.
$zhcn.Range("P3").Value = "Handback file name: bkumuoun.0xo is different with handoff file name: e50d8f56-22e6-4e33-9e12-6675c1567d4b.a4026c43e9a87e1896a693d608773986baebce5b.zh-cn."
$dede.Range("P3").Value = "Handback file name: bkumuoun.0xo is different with handoff file name: e50d8f56-22e6-4e33-9e12-6675c1567d4b.a4026c43e9a87e1896a693d608773986baebce5b.de-de."

# Widen the Error Detail column to fit the new, longer messages.
# (ColumnWidth is specified in characters; Excel pads/rounds this to a
# stored width in 1/256-character units when writing the file, so 39.1667
# is the character-width input that round-trips to a stored width of 40.)
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667
